$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: num_customers (C34) 61 -> 62, retention_rate (E34) recalculated as C34/D34
$ws.Range("C34").Value = 62
$ws.Range("E34").Value = 62/2256

# Row 37: num_customers (C37) and cohort_size (D37) 599 -> 604
$ws.Range("C37").Value = 604
$ws.Range("D37").Value = 604
